$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.418.84'
$ws.Range('E2').Value = '  -2.98%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.258.11'
$ws.Range('E3').Value = '  -5.72%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '589.98'
$ws.Range('E5').Value = '  -3.35%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '150.25'
$ws.Range('E6').Value = '  -10.38%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.250.60'
$ws.Range('E8').Value = '  -5.76%  '

$ws.Range('E9').Value = '  -8.39%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').Value = '  -10.33%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.65'
$ws.Range('E11').Value = '  -5.61%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.504'
$ws.Range('E12').Value = '  -10.57%  '

$ws.Range('E13').Value = '  -8.37%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '38.30'
$ws.Range('E14').Value = '  -13.63%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.781.54'
$ws.Range('E15').Value = '  -5.90%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.480.58'
$ws.Range('E16').Value = '  -3.07%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.260.23'
$ws.Range('E17').Value = '  -5.82%  '

$ws.Range('E18').Value = '  -5.38%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '530.05'
$ws.Range('E19').Value = '  -8.40%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.11'
$ws.Range('E20').Value = '  -13.05%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.95'
$ws.Range('E21').Value = '  -12.97%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.755'
$ws.Range('E22').Value = '  -10.83%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.87'
$ws.Range('E23').Value = '  -11.83%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.43'
$ws.Range('E24').Value = '  -10.80%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.49'
$ws.Range('E25').Value = '  -11.28%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.11%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.22'
$ws.Range('E27').Value = '  -11.17%  '

$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.15'
$ws.Range('E28').Value = '  -11.74%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.01'
$ws.Range('E29').Value = '  -7.06%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '29.09'
$ws.Range('E30').Value = '  -11.32%  '

$ws.Range('E31').Value = '  -5.09%  '

$ws.Range('E32').Value = '  -5.03%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.60'
$ws.Range('E33').Value = '  -16.01%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.71'
$ws.Range('E34').Value = '  -13.15%  '

$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.13%  '

$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '511.39'
$ws.Range('E36').Value = '  -12.76%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0439'
$ws.Range('E37').Value = '  -7.24%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '53.13'
$ws.Range('E38').Value = '  -5.25%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0854'
$ws.Range('E39').Value = '  -10.68%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.94'
$ws.Range('E40').Value = '  -15.28%  '

$ws.Range('E41').Value = '  -10.59%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.74'
$ws.Range('E42').Value = '  -12.86%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.930.05'
$ws.Range('E43').Value = '  -9.79%  '

$ws.Range('E44').Value = '  -10.16%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0₃0588'
$ws.Range('E45').Value = '  -15.51%  '

$ws.Range('E46').Value = '  -8.80%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '26.67'
$ws.Range('E47').Value = '  -13.72%  '

$ws.Range('E48').Value = '  -0.07%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.32'
$ws.Range('E49').Value = '  -16.64%  '

$ws.Range('E50').Value = '  -10.09%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '123.54'
$ws.Range('E51').Value = '  -7.55%  '

